$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that changes from 45178 (2023-09-09)
# to 45179 (2023-09-10) for every data row (rows 2 through 261).
for ($row = 2; $row -le 261; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
